$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: resize/reposition the document window to match the author's
# recorded window geometry (xWindow/yWindow/windowWidth/windowHeight).
try {
    $win = $excel.ActiveWindow
    $win.Left = 13480
    $win.Top = 0
    $win.Width = 37560
    $win.Height = 26000
} catch {
}

# Update the Phytoplankton Functional Type (PFT) labels in column A
# (rows 3,4,7 are unchanged; rows 5,6,8,9 get reworded labels)
$ws.Range("A5").Value = "Prymensiophytes (chromophytes and nanoflagellates)"
$ws.Range("A6").Value = "Pelagophytes (chromophytes and nanoflagellates)"
$ws.Range("A8").Value = "Green algae (green flagellates and prochlorophytes)"
$ws.Range("A9").Value = "Prokaryotes (cyanobacteria and prochlorophytes)"

# Widen column A so the longer labels are readable (was 27.6640625 -> 38.5)
$ws.Columns("A").ColumnWidth = 37.67

# Move the active selection from A11 to A6
$ws.Range("A6").Select()
